$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "!" marker column (old column E) entirely; this shifts the
# old F/G date columns left into E/F for rows 10-12.
$ws.Columns("E").Delete()

# Add the two new header cells for the additional mention columns.
$ws.Range("E1").Value = "additional_mention1"
$ws.Range("F1").Value = "additional_mention_2"

# Leave the selection where the editor last left it.
[void]$ws.Range("F6").Select()
